$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "API" section of the sheet (rows 18-29) is being restructured:
#  - The "Leads" sub-section (rows 19-23) is replaced by an "Accounts" sub-section,
#    which also gains a new "Get Accounts List" row.
#  - The "Contacts" sub-section gains Type values (Smoke/Regres) and a new
#    "Get Contacts List" row.
#  - The trailing "e2e" rows shift down to make room.
# Clear out the old rows 18-29 (A:E) first, then rewrite rows 18-31 from scratch.
$ws.Range("A18:E29").ClearContents()

# Row 18 - API section header (unchanged)
$ws.Range("A18").Value = "API"

# Rows 19-24 - Accounts (API)
$ws.Range("A19").Value = "Accounts"
$ws.Range("B19").Value = "Create Account"
$ws.Range("E19").Value = "Smoke"

$ws.Range("B20").Value = "Read Account"
$ws.Range("E20").Value = "Smoke"

$ws.Range("B21").Value = "Edit Account"
$ws.Range("E21").Value = "Regres"

$ws.Range("B22").Value = "Delete Account"
$ws.Range("E22").Value = "Regres"

$ws.Range("B23").Value = "Create account with empty fields"
$ws.Range("E23").Value = "Regres"

$ws.Range("B24").Value = "Get Accounts List"
$ws.Range("E24").Value = "Smoke"

# Rows 25-29 - Contacts (API)
$ws.Range("A25").Value = "Contacts"
$ws.Range("B25").Value = "Create Contact"
$ws.Range("E25").Value = "Smoke"

$ws.Range("B26").Value = "Read Contact"
$ws.Range("E26").Value = "Smoke"

$ws.Range("B27").Value = "Edit Contact"
$ws.Range("E27").Value = "Regres"

$ws.Range("B28").Value = "Delete Contact"
$ws.Range("E28").Value = "Regres"

$ws.Range("B29").Value = "Get Contacts List"
$ws.Range("E29").Value = "Smoke"

# Rows 30-31 - e2e (shifted down from 28-29)
$ws.Range("A30").Value = "e2e"
$ws.Range("B30").Value = "Create New Account"
$ws.Range("C30").Value = "One 3-step scenario"
$ws.Range("E30").Value = "e2e"

$ws.Range("B31").Value = "Create New Contact"
$ws.Range("C31").Value = "Use cucumber?"
$ws.Range("E31").Value = "e2e"

# Update the view: scrolled/selected cell moved to E22 in the saved file.
$ws.Range("E22").Select()
